$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.623.37"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "2.113.31"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.85"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4506"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.58"
$ws.Range("E9").Value = "  +0.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09036"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.38"
$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("D13").Value = "2.125.75"
$ws.Range("E13").Value = "  +1.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.782"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.061"
$ws.Range("E15").Value = "  +3.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.90"
$ws.Range("E16").Value = "  +1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001163"
$ws.Range("E17").Value = "  +2.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.013"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06696"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.35"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.349"
$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("D23").Value = "30.711.35"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("E24").Value = "  +3.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.376"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").Value = "2.371.94"
$ws.Range("E26").Value = "  +1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.39"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.34"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.81"
$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.377"
$ws.Range("E33").Value = "  +3.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.632"
$ws.Range("E34").Value = "  -2.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.940"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.30"
$ws.Range("E36").Value = "  -2.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.895"
$ws.Range("E37").Value = "  +5.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02651"
$ws.Range("E38").Value = "  +2.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06828"
$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2320"
$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.60"
$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6873"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.02"
$ws.Range("E44").Value = "  +6.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6425"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.313"
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000368"
$ws.Range("E47").Value = "  +12.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.702"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.253"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.94"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07303"
$ws.Range("E51").Value = "  +3.03%  "
